$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be auto-parsed as numbers
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply updated price (D) and volume/% change (E) values
$ws.Range('D2').Value = '62.721.39'
$ws.Range('E2').Value = '  -1.62%  '
$ws.Range('D3').Value = '3.030.37'
$ws.Range('E3').Value = '  -1.94%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '582.75'
$ws.Range('E5').Value = '  -1.78%  '
$ws.Range('D6').Value = '149.20'
$ws.Range('E6').Value = '  -4.55%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  -3.08%  '
$ws.Range('D9').Value = '3.030.05'
$ws.Range('E9').Value = '  -1.87%  '
$ws.Range('E10').Value = '  -3.75%  '
$ws.Range('D11').Value = '5.66'
$ws.Range('E11').Value = '  -3.07%  '
$ws.Range('E12').Value = '  -2.41%  '
$ws.Range('E13').Value = '  -3.98%  '
$ws.Range('D14').Value = '35.35'
$ws.Range('E14').Value = '  -5.66%  '
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').Value = '3.533.89'
$ws.Range('E16').Value = '  -1.92%  '
$ws.Range('D17').Value = '7.06'
$ws.Range('E17').Value = '  -1.73%  '
$ws.Range('D18').Value = '62.710.84'
$ws.Range('D19').Value = '3.029.16'
$ws.Range('E19').Value = '  -2.03%  '
$ws.Range('D20').Value = '468.04'
$ws.Range('E20').Value = '  -2.50%  '
$ws.Range('D21').Value = '14.07'
$ws.Range('E21').Value = '  -3.70%  '
$ws.Range('E22').Value = '  -2.89%  '
$ws.Range('D23').Value = '7.41'
$ws.Range('E23').Value = '  -2.15%  '
$ws.Range('D24').Value = '2.39'
$ws.Range('E24').Value = '  -1.19%  '
$ws.Range('D25').Value = '81.07'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('E26').Value = '  -3.46%  '
$ws.Range('D27').Value = '10.52'
$ws.Range('E27').Value = '  +2.26%  '
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '7.25'
$ws.Range('E30').Value = '  -3.48%  '
$ws.Range('D31').Value = '2.63'
$ws.Range('E31').Value = '  -2.16%  '
$ws.Range('E32').Value = '  -1.33%  '
$ws.Range('D33').Value = '27.57'
$ws.Range('E33').Value = '  +0.91%  '
$ws.Range('E34').Value = '  -5.06%  '
$ws.Range('E35').Value = '  -1.10%  '
$ws.Range('D36').Value = '0.0₃0797'
$ws.Range('E36').Value = '  -6.77%  '
$ws.Range('E37').Value = '  -4.82%  '
$ws.Range('E38').Value = '  -2.97%  '
$ws.Range('D39').Value = '50.27'
$ws.Range('E39').Value = '  -1.15%  '
$ws.Range('E40').Value = '  -15.01%  '
$ws.Range('D42').Value = '421.93'
$ws.Range('E42').Value = '  -5.68%  '
$ws.Range('D43').Value = '0.283'
$ws.Range('E43').Value = '  -1.91%  '
$ws.Range('E44').Value = '  +0.80%  '
$ws.Range('D45').Value = '2.796.50'
$ws.Range('E45').Value = '  -1.22%  '
$ws.Range('E46').Value = '  -2.02%  '
$ws.Range('D47').Value = '38.02'
$ws.Range('E47').Value = '  -9.96%  '
$ws.Range('D48').Value = '129.84'
$ws.Range('E48').Value = '  -1.43%  '
$ws.Range('D50').Value = '24.54'
$ws.Range('E50').Value = '  -4.36%  '
$ws.Range('E51').Value = '  -1.61%  '
